# Natmi following Dr Hou advice
# A new sending/target cluster "M2" was added to the LR-pair analysis
# (Slit2-Robo4), which changes the aggregate statistics for all rows
# and appends two new rows (M2 -> ECs, M2 -> sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data block, A2:T9 (8 data rows, 20 columns).
# Sending-cluster order is now ECs, FAPs, M2, sCs (M2 newly inserted
# between FAPs and sCs); each sending cluster has a row for each of the
# two target clusters (ECs, sCs).
$data = @(
    @("ECs","Slit2","Robo4","ECs",3,1,0.09783033333333334,0.293491,0.01145526398792098,0.01145526398792098,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,4.826112229926777,43.435010069341,0.01144374062914415,0.01144374062914415),
    @("ECs","Slit2","Robo4","sCs",3,1,0.09783033333333334,0.293491,0.01145526398792098,0.01145526398792098,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,0.004859689198222222,0.043737202784,0.00001152335877683316,0.00001152335877683316),
    @("FAPs","Slit2","Robo4","ECs",3,1,8.296819666666666,24.890459,0.9715009272022777,0.9715009272022777,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,409.2941473107898,3683.647325797108,0.970523651274985,0.970523651274985),
    @("FAPs","Slit2","Robo4","sCs",3,1,8.296819666666666,24.890459,0.9715009272022777,0.9715009272022777,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,0.4121417513351111,3.709275762015999,0.00097727592729268,0.00097727592729268),
    @("M2","Slit2","Robo4","ECs",1,0.3333333333333333,0.002823666666666667,0.008470999999999999,0.0003306320849418846,0.0003306320849418846,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,0.1392955719245556,1.253660147321,0.0003302994874441809,0.0003302994874441809),
    @("M2","Slit2","Robo4","sCs",1,0.3333333333333333,0.002823666666666667,0.008470999999999999,0.0003306320849418846,0.0003306320849418846,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,0.0001402647004444444,0.001262382304,0.0000003325974977036901,0.0000003325974977036901),
    @("sCs","Slit2","Robo4","ECs",3,1,0.142734,0.428202,0.0167131767248595,0.0167131767248595,3,1,49.33145033333333,147.994351,0.9989940555897285,0.9989940555897285,7.041275231877999,63.37147708690199,0.01669636419815525,0.01669636419815525),
    @("sCs","Slit2","Robo4","sCs",3,1,0.142734,0.428202,0.0167131767248595,0.0167131767248595,3,1,0.04967466666666667,0.149024,0.001005944410271468,0.001005944410271468,0.007090263871999999,0.06381237484799999,0.00001681252670425163,0.00001681252670425163)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
